# use 3rd quartile instead of mean
# Update Number_of_Inclusions (col B) and the recomputed
# Number_of_Inclusions_per_Nucleus (col D = B / C) for the affected rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 3;  B = 0;  D = 0 },
    @{ Row = 11; B = 6;  D = 0.2142857142857143 },
    @{ Row = 14; B = 5;  D = 0.1724137931034483 },
    @{ Row = 15; B = 0;  D = 0 },
    @{ Row = 22; B = 6;  D = 0.2857142857142857 },
    @{ Row = 23; B = 28; D = 1.473684210526316 },
    @{ Row = 24; B = 0;  D = 0 },
    @{ Row = 31; B = 1;  D = 0.05 },
    @{ Row = 37; B = 4;  D = 0.2666666666666667 },
    @{ Row = 40; B = 1;  D = 0.05882352941176471 },
    @{ Row = 47; B = 0;  D = 0 },
    @{ Row = 56; B = 2;  D = 0.07407407407407407 }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 2).Value = $u.B
    $ws.Cells.Item($u.Row, 4).Value = $u.D
}
